$d = $word.ActiveDocument

# 1. Replace the ID placeholder text in the first paragraph.
$d.Content.Find.Execute("**ID__AFFARS_5327_topic_2__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_SUBPART_5327_2__ID**", 2)

# 2. Remove the trailing space run that used to follow the ID placeholder
#    (it sits right before the paragraph mark of paragraph 1).
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$trailingSpace = $d.Range($r1.End - 2, $r1.End - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# 3. Update paragraph formatting: indent and paragraph border (space-only, no line).
$p1 = $d.Paragraphs(1)
$p1.Format.LeftIndent = 225 / 20

$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
